# Apply updated crypto price/volume data as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.833.53"
$ws.Range("D3").Value = "2.619.24"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("D5").Value = "'604.20"
$ws.Range("E5").Value = "  +1.62%  "
$ws.Range("D6").Value = "'154.41"
$ws.Range("E6").Value = "  -0.49%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +1.64%  "
$ws.Range("D9").Value = "2.615.42"
$ws.Range("E9").Value = "  +0.83%  "
$ws.Range("D10").Value = "'0.128"
$ws.Range("E10").Value = "  +11.50%  "
$ws.Range("E13").Value = "  -0.99%  "
$ws.Range("D14").Value = "'27.75"
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("D15").Value = "'0.0000188"
$ws.Range("E15").Value = "  +3.91%  "
$ws.Range("D16").Value = "3.094.40"
$ws.Range("E16").Value = "  +1.04%  "
$ws.Range("D17").Value = "67.732.63"
$ws.Range("E17").Value = "  +1.32%  "
$ws.Range("D18").Value = "2.623.22"
$ws.Range("E18").Value = "  +0.95%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "'367.31"
$ws.Range("E19").Value = "  +3.09%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'11.18"
$ws.Range("E20").Value = "  -1.02%  "
$ws.Range("D21").Value = "'7.66"
$ws.Range("E21").Value = "  -1.56%  "
$ws.Range("D22").Value = "'4.31"
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("D23").Value = "'2.05"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("E25").Value = "  +4.38%  "
$ws.Range("D26").Value = "'9.87"
$ws.Range("E26").Value = "  -5.71%  "
$ws.Range("D27").Value = "'0.0000105"
$ws.Range("E27").Value = "  +1.06%  "
$ws.Range("D28").Value = "2.744.76"
$ws.Range("E28").Value = "  +1.22%  "
$ws.Range("D29").Value = "'580.43"
$ws.Range("E29").Value = "  -3.01%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  -2.56%  "
$ws.Range("D32").Value = "'7.92"
$ws.Range("E32").Value = "  -1.76%  "
$ws.Range("D33").Value = "'1.87"
$ws.Range("E33").Value = "  +1.04%  "
$ws.Range("E34").Value = "  -1.77%  "
$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("D36").Value = "'1.53"
$ws.Range("E36").Value = "  -3.38%  "
$ws.Range("D37").Value = "'4.94"
$ws.Range("E37").Value = "  -1.48%  "
$ws.Range("D38").Value = "'158.07"
$ws.Range("E38").Value = "  +2.76%  "
$ws.Range("D39").Value = "'19.42"
$ws.Range("E39").Value = "  +1.00%  "
$ws.Range("D40").Value = "'0.370"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("D41").Value = "'5.37"
$ws.Range("E41").Value = "  -2.19%  "
$ws.Range("D43").Value = "'2.62"
$ws.Range("E43").Value = "  -1.18%  "
$ws.Range("D44").Value = "'41.20"
$ws.Range("E44").Value = "  -0.77%  "
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").Value = "'16.43"
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("D47").Value = "'156.90"
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("D48").Value = "0.0₆0287"
$ws.Range("E48").Value = "  -7.60%  "
$ws.Range("D49").Value = "'3.76"
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("D50").Value = "'20.98"
$ws.Range("E50").Value = "  -2.21%  "
$ws.Range("B51").Value = "Hedera"
$ws.Range("C51").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D51").Value = "'0.0540"
$ws.Range("E51").Value = "  -3.44%  "

Write-Output "Applied all cell updates"
